$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the Time Log entry for row 91 (previously blank placeholder row)
$ws.Range("A91").Value = 41930
$ws.Range("B91").Value = 0.79861111111111116
$ws.Range("C91").Value = 0.81388888888888899
$ws.Range("D91").Value = 0
$ws.Range("F91").Value = "Coding"

# Move the active selection to C92, matching the recorded cursor position
$ws.Range("C92").Select()

$wb.Save()
